$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 642.8
$ws.Range("I5").Value = 642.8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 642.8
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -527.8

$ws.Range("H32").Value = 3393.6
$ws.Range("I32").Value = 1191.6666
$ws.Range("K32").Value = 1191.6666
$ws.Range("M32").Value = -865.6666

$ws.Range("H96").Value = 418.42856
$ws.Range("I96").Value = 418.42856
$ws.Range("K96").Value = 1255.28568
$ws.Range("M96").Value = 117.71432

$ws.Range("H132").Value = 1987.5834
$ws.Range("I132").Value = 1972.7593
$ws.Range("J132").Value = 2121
$ws.Range("K132").Value = 5918.2779
$ws.Range("L132").Value = 6363
$ws.Range("M132").Value = -3388.2779
$ws.Range("N132").Value = -11423

$ws.Range("H138").Value = 3441.0483
$ws.Range("I138").Value = 2434.2307
$ws.Range("K138").Value = 7302.6921
$ws.Range("M138").Value = -2162.6921

$ws.Range("H141").Value = 19833.834
$ws.Range("I141").Value = 35680.668
$ws.Range("J141").Value = 3987
$ws.Range("K141").Value = 107042.004
$ws.Range("L141").Value = 11961
$ws.Range("M141").Value = -101862.004
$ws.Range("N141").Value = -22321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3143939
$ws.Range("I2").Value = 5143442
$ws.Range("J2").Value = 1863.4286
$ws.Range("K2").Value = 5143442
$ws.Range("L2").Value = 1863.4286
$ws.Range("M2").Value = -5143329
$ws.Range("N2").Value = -2089.4286

$ws.Range("H32").Value = 7768.9785
$ws.Range("I32").Value = 6382.6943
$ws.Range("K32").Value = 6382.6943
$ws.Range("M32").Value = -6095.6943

$ws.Range("H61").Value = 4500.7417
$ws.Range("I61").Value = 4093.6296
$ws.Range("J61").Value = 7248.75
$ws.Range("K61").Value = 4093.6296
$ws.Range("L61").Value = 7248.75
$ws.Range("M61").Value = -3881.6296
$ws.Range("N61").Value = -7672.75

$ws.Range("H74").Value = 21209.809
$ws.Range("I74").Value = 1738.5135
$ws.Range("K74").Value = 1738.5135
$ws.Range("M74").Value = -864.5135

$ws.Range("H77").Value = 21209.809
$ws.Range("I77").Value = 1738.5135
$ws.Range("K77").Value = 8692.567500000001
$ws.Range("M77").Value = -4324.567500000001

$ws.Range("H97").Value = 1446756.4
$ws.Range("I97").Value = 2314531.5
$ws.Range("K97").Value = 2314531.5
$ws.Range("M97").Value = -2314035.5

$ws.Range("H116").Value = 3143939
$ws.Range("I116").Value = 5143442
$ws.Range("J116").Value = 1863.4286
$ws.Range("K116").Value = 5143442
$ws.Range("L116").Value = 1863.4286
$ws.Range("M116").Value = -5141148
$ws.Range("N116").Value = -6451.4286

$ws.Range("H132").Value = 36609.652
$ws.Range("I132").Value = 1685.55
$ws.Range("J132").Value = 153023.33
$ws.Range("K132").Value = 5056.65
$ws.Range("L132").Value = 459069.99
$ws.Range("M132").Value = -2526.65
$ws.Range("N132").Value = -464129.99

$ws.Range("H136").Value = 4500.7417
$ws.Range("I136").Value = 4093.6296
$ws.Range("J136").Value = 7248.75
$ws.Range("K136").Value = 12280.8888
$ws.Range("L136").Value = 21746.25
$ws.Range("M136").Value = -9730.888800000001
$ws.Range("N136").Value = -26846.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3143939
$ws.Range("I3").Value = 5143442
$ws.Range("J3").Value = 1863.4286
$ws.Range("K3").Value = 5143442
$ws.Range("L3").Value = 1863.4286
$ws.Range("M3").Value = -5143328
$ws.Range("N3").Value = -2091.4286

$ws.Range("H81").Value = 23374.75
$ws.Range("J81").Value = 23374.75
$ws.Range("L81").Value = 23374.75
$ws.Range("N81").Value = -25496.75

$ws.Range("H84").Value = 23374.75
$ws.Range("J84").Value = 23374.75
$ws.Range("L84").Value = 70124.25
$ws.Range("N84").Value = -80732.25

$ws.Range("H94").Value = 4812941
$ws.Range("I94").Value = 5954594
$ws.Range("K94").Value = 5954594
$ws.Range("M94").Value = -5954143

$ws.Range("H99").Value = 7520340.5
$ws.Range("I99").Value = 11905844
$ws.Range("K99").Value = 11905844
$ws.Range("M99").Value = -11904346

$ws.Range("H105").Value = 3908961.8
$ws.Range("I105").Value = 5211075
$ws.Range("K105").Value = 5211075
$ws.Range("M105").Value = -5209328

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25531
$ws.Range("J31").Value = 74613.08
$ws.Range("L31").Value = 74613.08
$ws.Range("N31").Value = -75203.08

$ws.Range("H34").Value = 25531
$ws.Range("J34").Value = 74613.08
$ws.Range("L34").Value = 74613.08
$ws.Range("N34").Value = -75017.08

$ws.Range("H86").Value = 6297.875
$ws.Range("I86").Value = 5040.9546
$ws.Range("J86").Value = 9063.1
$ws.Range("K86").Value = 5040.9546
$ws.Range("L86").Value = 9063.1
$ws.Range("M86").Value = -3917.9546
$ws.Range("N86").Value = -11309.1

$ws.Range("H89").Value = 6297.875
$ws.Range("I89").Value = 5040.9546
$ws.Range("J89").Value = 9063.1
$ws.Range("K89").Value = 25204.773
$ws.Range("L89").Value = 45315.5
$ws.Range("M89").Value = -19588.773
$ws.Range("N89").Value = -56547.5

$ws.Range("H134").Value = 1664.8334
$ws.Range("I134").Value = 1092.2858
$ws.Range("J134").Value = 4527.5713
$ws.Range("K134").Value = 3276.8574
$ws.Range("L134").Value = 13582.7139
$ws.Range("M134").Value = -741.8574000000003
$ws.Range("N134").Value = -18652.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17365850
$ws.Range("I131").Value = 20834152
$ws.Range("J131").Value = 16672191
$ws.Range("K131").Value = 62502456
$ws.Range("L131").Value = 50016573
$ws.Range("M131").Value = -62497416
$ws.Range("N131").Value = -50026653

$ws.Range("H137").Value = 5833.12
$ws.Range("J137").Value = 6270
$ws.Range("L137").Value = 18810
$ws.Range("N137").Value = -29010

$ws.Range("H138").Value = 5314.2666
$ws.Range("I138").Value = 2624.8333
$ws.Range("J138").Value = 7107.222
$ws.Range("K138").Value = 7874.499899999999
$ws.Range("L138").Value = 21321.666
$ws.Range("M138").Value = -2734.499899999999
$ws.Range("N138").Value = -31601.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 744619.1
$ws.Range("J97").Value = 582.63635
$ws.Range("L97").Value = 582.63635
$ws.Range("N97").Value = -1574.63635

$ws.Range("H113").Value = 7363071.5
$ws.Range("I113").Value = 11042815
$ws.Range("K113").Value = 11042815
$ws.Range("M113").Value = -11040645

$ws.Range("H132").Value = 3438.4092
$ws.Range("I132").Value = 3152.0908
$ws.Range("K132").Value = 9456.2724
$ws.Range("M132").Value = -6926.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7913.1763
$ws.Range("I40").Value = 4074.111
$ws.Range("K40").Value = 4074.111
$ws.Range("M40").Value = -3938.111

$ws.Range("H46").Value = 4394.6924
$ws.Range("I46").Value = 878.4286
$ws.Range("J46").Value = 8497
$ws.Range("K46").Value = 878.4286
$ws.Range("L46").Value = 8497
$ws.Range("M46").Value = -690.4286
$ws.Range("N46").Value = -8873

$ws.Range("H63").Value = 192538.5
$ws.Range("J63").Value = 85000
$ws.Range("L63").Value = 85000
$ws.Range("N63").Value = -86498

$ws.Range("H66").Value = 192538.5
$ws.Range("J66").Value = 85000
$ws.Range("L66").Value = 255000
$ws.Range("N66").Value = -262488

$ws.Range("H68").Value = 1800
$ws.Range("I68").Value = 1950
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -1201
$ws.Range("N68").Value = -2998

$ws.Range("H71").Value = 1800
$ws.Range("I71").Value = 1950
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 9750
$ws.Range("L71").Value = 7500
$ws.Range("M71").Value = -6006
$ws.Range("N71").Value = -14988
